# Update Excel file with latest predictions (24-01-2025 matches)

$wb = $excel.ActiveWorkbook

function Set-RowValues {
    param($sheet, $row, $values)
    $col = 1
    foreach ($v in $values) {
        $sheet.Cells.Item($row, $col).Value = $v
        $col = $col + 1
    }
}

# --- Sheet "Home win": append two new rows (4 and 5) ---
$wsHome = $wb.Worksheets.Item("Home win")
Set-RowValues $wsHome 4 @("24-01-2025 19:00", "NETHERLANDS", "EERSTE DIVISIE", "Vitesse - Dordrecht", 73.3, 2.7)
Set-RowValues $wsHome 5 @("24-01-2025 23:30", "WORLD", "SUDAMERICANO U20", "Brazil U20 - Argentina U20", 70, 2.45)

# --- Sheet "Draw": append one new row (3) ---
$wsDraw = $wb.Worksheets.Item("Draw")
Set-RowValues $wsDraw 3 @("24-01-2025 13:30", "INDIA", "I-LEAGUE", "Gokulam - Inter Kashi", 66.7, 3.6)

# --- Sheet "Btts": replace rows 2-5 with new matches, append row 6 ---
$wsBtts = $wb.Worksheets.Item("Btts")
Set-RowValues $wsBtts 2 @("24-01-2025 19:30", "GERMANY", "BUNDESLIGA", "VfL Wolfsburg - Holstein Kiel", 76, 1.75)
Set-RowValues $wsBtts 3 @("24-01-2025 17:00", "CROATIA", "HNL", "NK Osijek - Sibenik", 76.7, 2)
Set-RowValues $wsBtts 4 @("24-01-2025 10:45", "ISRAEL", "LIGA ALEF", "Hapoel Marmorek - Agudat Sport Ashdod", 76.7, 1.77)
Set-RowValues $wsBtts 5 @("24-01-2025 19:30", "ITALY", "SERIE C - GIRONE B", "Lucchese - Ascoli", 84, 1.83)
Set-RowValues $wsBtts 6 @("24-01-2025 17:00", "TURKEY", "SÜPER LIG", "Samsunspor - Gazişehir Gaziantep", 76.7, 1.73)

# --- Sheet "Over_Under": append six new rows (5-10) ---
$wsOU = $wb.Worksheets.Item("Over_Under")
Set-RowValues $wsOU 5 @("24-01-2025 19:30", "GERMANY", "BUNDESLIGA", "VfL Wolfsburg - Holstein Kiel", 86.7, 1.57, 60, 2.5)
Set-RowValues $wsOU 6 @("24-01-2025 17:00", "CROATIA", "HNL", "NK Osijek - Sibenik", 80, 1.75, 50, 2.9)
Set-RowValues $wsOU 7 @("24-01-2025 19:00", "NETHERLANDS", "EERSTE DIVISIE", "FC Eindhoven - Cambuur", 80, 1.73, 55, 2.8)
Set-RowValues $wsOU 8 @("24-01-2025 11:10", "ISRAEL", "LIGA ALEF", "Maccabi Ahi Nazareth - Hapoel Beit Shean", 80, 1.8, 33.3, 2.95)
Set-RowValues $wsOU 9 @("24-01-2025 19:15", "SWITZERLAND", "CHALLENGE LEAGUE", "Étoile Carouge - FC WIL 1900", 70, 1.73, 60, 2.88)
Set-RowValues $wsOU 10 @("24-01-2025 17:00", "TURKEY", "SÜPER LIG", "Samsunspor - Gazişehir Gaziantep", 80, 1.75, 40, 2.88)

# --- Sheet "Away Win": replace row 2 with new match ---
$wsAway = $wb.Worksheets.Item("Away Win")
Set-RowValues $wsAway 2 @("24-01-2025 19:00", "FRANCE", "LIGUE 2", "Martigues - Amiens", 70, 1.91)
